$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '30.174.50'
$ws.Range('E2').Value = '  +1.10%  '

# Row 3
$ws.Range('D3').Value = '1.895.07'
$ws.Range('E3').Value = '  +0.35%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.12%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7402'
$ws.Range('E5').Value = '  -0.95%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '243.23'
$ws.Range('E6').Value = '  +0.36%  '

# Row 7
$ws.Range('E7').Value = '  +0.12%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3174'
$ws.Range('E8').Value = '  +1.81%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07224'
$ws.Range('E9').Value = '  +1.38%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.97'
$ws.Range('E10').Value = '  -1.32%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08359'
$ws.Range('E11').Value = '  -1.55%  '

# Row 12
$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7616'
$ws.Range('E12').Value = '  +0.21%  '

# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.937.33'
$ws.Range('E13').Value = '  +4.61%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.463'
$ws.Range('E14').Value = '  +1.93%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.13'
$ws.Range('E15').Value = '  -0.27%  '

# Row 16
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '30.280.48'
$ws.Range('E16').Value = '  +1.43%  '

# Row 17
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.168'
$ws.Range('E17').Value = '  -0.07%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '250.95'
$ws.Range('E18').Value = '  +3.04%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.67'
$ws.Range('E19').Value = '  -0.22%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007917'
$ws.Range('E20').Value = '  +1.52%  '

# Row 21
$ws.Range('D21').Value = '2.185.49'
$ws.Range('E21').Value = '  +2.17%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.004'
$ws.Range('E22').Value = '  +0.43%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.990'
$ws.Range('E23').Value = '  -0.10%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.002'
$ws.Range('E24').Value = '  -0.06%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1582'
$ws.Range('E25').Value = '  -0.67%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.326'
$ws.Range('E26').Value = '  -0.25%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '164.73'
$ws.Range('E27').Value = '  +1.44%  '

# Row 28
$ws.Range('E28').Value = '  +0.33%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.070'
$ws.Range('E29').Value = '  +2.08%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.486'
$ws.Range('E30').Value = '  -0.70%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.597'
$ws.Range('E31').Value = '  +2.26%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.537'
$ws.Range('E32').Value = '  +0.33%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.213'
$ws.Range('E33').Value = '  +2.39%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05383'
$ws.Range('E34').Value = '  -0.42%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.258'
$ws.Range('E35').Value = '  +1.59%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7864'
$ws.Range('E36').Value = '  +5.41%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.007'
$ws.Range('E37').Value = '  +0.48%  '

# Row 38
$ws.Range('E38').Value = '  +0.82%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01973'
$ws.Range('E39').Value = '  +1.83%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.765'
$ws.Range('E40').Value = '  -0.31%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4582'
$ws.Range('E41').Value = '  +2.93%  '

# Row 42
$ws.Range('D42').Value = '1.102.03'
$ws.Range('E42').Value = '  +0.99%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.092'
$ws.Range('E43').Value = '  +0.38%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '72.99'
$ws.Range('E44').Value = '  +0.64%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8734'
$ws.Range('E45').Value = '  +2.20%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '104.73'
$ws.Range('E46').Value = '  +2.33%  '

# Row 47
$ws.Range('E47').Value = '  +0.23%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.870'
$ws.Range('E48').Value = '  +0.53%  '

# Row 49
$ws.Range('E49').Value = '  -0.98%  '

# Row 50
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.085.04'
$ws.Range('E50').Value = '  +1.83%  '

# Row 51
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.636'
$ws.Range('E51').Value = '  -0.96%  '
